# Rename the "Assessment" header (column J) to "Eligibility".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J1").Value = "Eligibility"

# The header row now wraps onto more lines, so its row height grows.
$ws.Rows.Item(1).RowHeight = 59.5

# Leave the cursor on D9, matching the saved selection state.
$ws.Range("D9").Select() | Out-Null
